$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F12").Value = 152
$ws1.Range("F13").Value = 134
$ws1.Range("F19").Value = 410
$ws1.Range("F20").Value = 987
$ws1.Range("F21").Value = 1594

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 152
$ws4.Range("F13").Value = 134
$ws4.Range("F19").Value = 410
$ws4.Range("F20").Value = 987
$ws4.Range("F21").Value = 1594
